# Auto-generated Excel COM-interop script
# Applies a scheduled market-price refresh to the Leve profit sheets
# (ALC, ARM, BSM, CUL, GSM, LTW, WVR) per the scraped commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 51: A Bile Business
$ws.Range("H51").Value = 3387.5715
$ws.Range("J51").Value = 3356.6667
$ws.Range("L51").Value = 3356.6667
$ws.Range("N51").Value = -4324.6667

# ALC row 98: The Dotted Line
$ws.Range("H98").Value = 803.6316
$ws.Range("I98").Value = 825.8125
$ws.Range("J98").Value = 685.3333
$ws.Range("K98").Value = 825.8125
$ws.Range("L98").Value = 685.3333
$ws.Range("M98").Value = 672.1875
$ws.Range("N98").Value = -3681.3333

# ALC row 105: Ultimate Official Strategy Guide
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

# ALC row 106: Making Your Mark
$ws.Range("H106").Value = 4675.5
$ws.Range("I106").Value = 3475.8333
$ws.Range("K106").Value = 3475.8333
$ws.Range("M106").Value = -2844.8333

# ALC row 116: Growing Up
$ws.Range("H116").Value = 3839.0789
$ws.Range("I116").Value = 4231.2383
$ws.Range("J116").Value = 3354.647
$ws.Range("K116").Value = 4231.2383
$ws.Range("L116").Value = 3354.647
$ws.Range("M116").Value = -789.2383
$ws.Range("N116").Value = -10238.647

# ALC row 122: Wishful Inking
$ws.Range("H122").Value = 803.6316
$ws.Range("I122").Value = 825.8125
$ws.Range("J122").Value = 685.3333
$ws.Range("K122").Value = 2477.4375
$ws.Range("L122").Value = 2055.9999
$ws.Range("M122").Value = -27.4375
$ws.Range("N122").Value = -6955.9999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 16: Greavous Losses
$ws.Range("H16").Value = 1203
$ws.Range("I16").Value = 937.3333
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 937.3333
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -650.3333
$ws.Range("N16").Value = -2574

# ARM row 19: Stadium Envy
$ws.Range("H19").Value = 32880.375
$ws.Range("I19").Value = 4254
$ws.Range("J19").Value = 61506.75
$ws.Range("K19").Value = 4254
$ws.Range("L19").Value = 61506.75
$ws.Range("M19").Value = -4025
$ws.Range("N19").Value = -61964.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 12: A Hit Job
$ws.Range("H12").Value = 375.9
$ws.Range("J12").Value = 980
$ws.Range("L12").Value = 980
$ws.Range("N12").Value = -1316

# BSM row 54: Get Me to the War on Time
$ws.Range("H54").Value = 5946.0625
$ws.Range("I54").Value = 2077.7
$ws.Range("K54").Value = 2077.7
$ws.Range("M54").Value = -1593.7

$ws = $wb.Worksheets.Item("CUL")
# CUL row 44: No More Dumpster Diving
$ws.Range("H44").Value = 759
$ws.Range("I44").Value = 138.25
$ws.Range("J44").Value = 2000.5
$ws.Range("K44").Value = 414.75
$ws.Range("L44").Value = 6001.5
$ws.Range("M44").Value = -16.75
$ws.Range("N44").Value = -6797.5

# CUL row 68: Such a Butter Face
$ws.Range("H68").Value = 523.9
$ws.Range("I68").Value = 467.8
$ws.Range("J68").Value = 580
$ws.Range("K68").Value = 1403.4
$ws.Range("L68").Value = 1740
$ws.Range("M68").Value = -592.4000000000001
$ws.Range("N68").Value = -3362

# CUL row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 523.9
$ws.Range("I71").Value = 467.8
$ws.Range("J71").Value = 580
$ws.Range("K71").Value = 4210.2
$ws.Range("L71").Value = 5220
$ws.Range("M71").Value = -154.1999999999998
$ws.Range("N71").Value = -13332

# CUL row 107: Slippery Service
$ws.Range("H107").Value = 368.72974
$ws.Range("I107").Value = 222.10527
$ws.Range("J107").Value = 523.5
$ws.Range("K107").Value = 666.3158099999999
$ws.Range("L107").Value = 1570.5
$ws.Range("M107").Value = 1253.68419
$ws.Range("N107").Value = -5410.5

# CUL row 122: Salt of the North
$ws.Range("H122").Value = 995.1429000000001
$ws.Range("I122").Value = 555.3333
$ws.Range("K122").Value = 4997.9997
$ws.Range("M122").Value = -2547.9997

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 1897.6364
$ws.Range("I131").Value = 2411.85
$ws.Range("J131").Value = 1717.2106
$ws.Range("K131").Value = 7235.549999999999
$ws.Range("L131").Value = 5151.6318
$ws.Range("M131").Value = -2195.549999999999
$ws.Range("N131").Value = -15231.6318

# CUL row 132: More Mezcal
$ws.Range("H132").Value = 11150.375
$ws.Range("J132").Value = 11199.833
$ws.Range("L132").Value = 100798.497
$ws.Range("N132").Value = -105858.497

# CUL row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 3032.2354
$ws.Range("I138").Value = 3032.2354
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 9096.706200000001
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3956.706200000001
$ws.Range("N138").ClearContents()

# CUL row 139: Najoothie
$ws.Range("H139").Value = 743.0769
$ws.Range("I139").Value = 586
$ws.Range("J139").Value = 1266.6666
$ws.Range("K139").Value = 1758
$ws.Range("L139").Value = 3799.9998
$ws.Range("M139").Value = 3382
$ws.Range("N139").Value = -14079.9998

# CUL row 141: Ocean Explosion
$ws.Range("H141").Value = 3898.9443
$ws.Range("I141").Value = 3974.0908
$ws.Range("K141").Value = 11922.2724
$ws.Range("M141").Value = -6742.2724

$ws = $wb.Worksheets.Item("GSM")
# GSM row 13: A Needle Is a Small Sword
$ws.Range("H13").Value = 745.375
$ws.Range("I13").Value = 331.4
$ws.Range("J13").Value = 1435.3334
$ws.Range("K13").Value = 331.4
$ws.Range("L13").Value = 1435.3334
$ws.Range("M13").Value = -192.4
$ws.Range("N13").Value = -1713.3334

# GSM row 14: All That Glitters
$ws.Range("H14").Value = 62970.938
$ws.Range("I14").Value = 77433.08
$ws.Range("J14").Value = 301.66666
$ws.Range("K14").Value = 77433.08
$ws.Range("L14").Value = 301.66666
$ws.Range("M14").Value = -77265.08
$ws.Range("N14").Value = -637.66666

# GSM row 17: Point of Honor
$ws.Range("H17").Value = 352
$ws.Range("J17").Value = 1000
$ws.Range("L17").Value = 1000
$ws.Range("N17").Value = -1336

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: Skin off Their Backs
$ws.Range("H22").Value = 731
$ws.Range("I22").Value = 400.1
$ws.Range("J22").Value = 1834
$ws.Range("K22").Value = 400.1
$ws.Range("L22").Value = 1834
$ws.Range("M22").Value = -105.1
$ws.Range("N22").Value = -2424

# LTW row 27: Fire and Hide
$ws.Range("H27").Value = 731
$ws.Range("I27").Value = 400.1
$ws.Range("J27").Value = 1834
$ws.Range("K27").Value = 400.1
$ws.Range("L27").Value = 1834
$ws.Range("M27").Value = -293.1
$ws.Range("N27").Value = -2048

# LTW row 41: The Hand that Bleeds
$ws.Range("H41").Value = 11545.667
$ws.Range("J41").Value = 11545.667
$ws.Range("L41").Value = 11545.667
$ws.Range("N41").Value = -12421.667

# LTW row 46: Supply Side Logic
$ws.Range("H46").Value = 1995.5
$ws.Range("I46").Value = 2326.6667
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 2326.6667
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -2138.6667
$ws.Range("N46").Value = -1378

$ws = $wb.Worksheets.Item("WVR")
# WVR row 8: The Adventurer's New Coat
$ws.Range("H8").Value = 751482.5
$ws.Range("I8").Value = 3000000
$ws.Range("J8").Value = 1976.6666
$ws.Range("K8").Value = 3000000
$ws.Range("L8").Value = 1976.6666
$ws.Range("M8").Value = -2999860
$ws.Range("N8").Value = -2256.6666

# WVR row 10: Just for Kecks
$ws.Range("H10").Value = 43192.4
$ws.Range("J10").Value = 43192.4
$ws.Range("L10").Value = 43192.4
$ws.Range("N10").Value = -43530.4

# WVR row 14: Hat in Hand
$ws.Range("H14").Value = 51504000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 51504000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 51504000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -51504336

# WVR row 19: Dirt Cheap
$ws.Range("H19").Value = 36904.25
$ws.Range("I19").Value = 3605
$ws.Range("J19").Value = 48004
$ws.Range("K19").Value = 3605
$ws.Range("L19").Value = 48004
$ws.Range("M19").Value = -3431
$ws.Range("N19").Value = -48352

# WVR row 31: Whatchoo Talking About
$ws.Range("H31").Value = 8750
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 8750
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 8750
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -9446

# WVR row 136: Weaving the Envelope
$ws.Range("H136").Value = 11490247
$ws.Range("I136").Value = 18201356
$ws.Range("K136").Value = 54604068
$ws.Range("M136").Value = -54601518
